# Update the demand-units source filename (refreshed version date)
# and the downstream demand-list indices range, per the updated
# delivery/demand extraction indices (only affects delivery data extraction).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B20: Demand Units File name -> new version dated 2025-11-30
$ws.Range("B20").Value = "cs3rpt2022_all_demand_units_v20251130.xlsx"

# D22: Demands List Block lower-right cell index shifted from O389 to O391
$ws.Range("D22").Value = "O391"

# Carry the D5:D11 "Lower Right Cell" formatting over to the new blank
# column E for those same rows (matches the widened used-range / E column
# that results from this edit in the authored workbook).
$ws.Range("D5:D11").Copy()
$ws.Range("E5:E11").PasteSpecial(-4122)

# Leave the active selection where the user made the last edit.
$ws.Range("D23").Select()
